# fix(publipostage): Correct status name
# Replace the "statut_label" value "bleu" -> "noir"
# and the "statut_name" value "pas de résultat ni de publication"
#   -> "pas de résultat postés ni publiés"
# for every data row in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldLabel = "bleu"
$newLabel = "noir"
$oldName  = "pas de résultat ni de publication"
$newName  = "pas de résultat postés ni publiés"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $bCell = $ws.Cells.Item($r, 2)
    if ($bCell.Value() -eq $oldLabel) {
        $bCell.Value = $newLabel
    }

    $cCell = $ws.Cells.Item($r, 3)
    if ($cCell.Value() -eq $oldName) {
        $cCell.Value = $newName
    }
}
